# edit.ps1 - apply "finished my part (1 and 5)" commit
#
# Summary of the change:
#   - Slide 1 title: "Crime" -> "Crimes in the US"
#   - Slides 2-6 "shift" their finished content forward by one slot and
#     slide 2/3 receive brand-new copy:
#       slide2: "Problem Statement and Background" -> "Problem Statement"
#               + new Problem Statement / Background paragraphs
#       slide3: "Methods Explored"        -> "Background" + new copy
#       slide4: "The Methods/Tools Used"  -> "Methods Explored" (old slide3 body)
#       slide5: "Results"                 -> "The Methods/Tools Used" (old slide4 body)
#       slide6: "Lessons Learned and Future Plans" -> "Results" (old slide5 body)
#   - Two new slides are appended at the end:
#       slide7: "Lessons Learned" (new bullet content)
#       slide8: "Future Plans" (new bullet content)
#
# Helper: write a list of paragraph strings into a TextRange without
# tripping the host's "multi-paragraph Text= assignment" sentence
# splitter - set the first paragraph via .Text, then grow the rest one
# paragraph at a time via InsertAfter so each ends up as a single run.
function Set-Paragraphs($textRange, [string[]]$paragraphs) {
    $textRange.Text = $paragraphs[0]
    for ($i = 1; $i -lt $paragraphs.Length; $i++) {
        $null = $textRange.InsertAfter("`r" + $paragraphs[$i])
    }
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Append two new slides at the end by duplicating the last existing
#    slide (keeps identical placeholder / layout / run formatting), then
#    overwrite their title + body text.
# ---------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$null = $lastSlide.Duplicate()
$null = $lastSlide.Duplicate()

$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Lessons Learned"
Set-Paragraphs $s7.Shapes.Item(2).TextFrame.TextRange @(
    "Analyzed crime trends to identify patterns and anomalies over time.",
    "Investigated the influence of poverty, education, and unemployment on crime rates.",
    "Evaluated the success of different policing strategies in crime reduction.",
    "Faced challenges in gathering accurate crime data.",
    "The benefits of integrating insights from various disciplines."
)

$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Future Plans"
Set-Paragraphs $s8.Shapes.Item(2).TextFrame.TextRange @(
    "Explore specific crime types and demographics in further detail.",
    "Suggest policy and law enforcement improvements.",
    "Collaborating with field experts for practical application and validation",
    "Utilize advanced data analysis techniques.",
    "Assessing the implications of findings on public policy and crime prevention"
)

# ---------------------------------------------------------------------
# 2. Slide 1: title text update.
# ---------------------------------------------------------------------
$p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange.Text = "Crimes in the US"

# ---------------------------------------------------------------------
# 3. Slide 2: "Problem Statement and Background" -> "Problem Statement"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Problem Statement"
Set-Paragraphs $s2.Shapes.Item(2).TextFrame.TextRange @(
    "“Analyzing Trends and Determinants of Crime Rates in the United States to Inform Policy and Prevention Strategies”",
    "In recent years, the United States has faced complex challenges related to crime. Understanding the trends, causes, and distribution of crime across different regions is crucial for developing effective policies and prevention strategies. This project aims to analyze various factors influencing crime rates, including economic, social, and environmental variables, to identify key drivers and potential areas for intervention."
)

# ---------------------------------------------------------------------
# 4. Slide 3: "Methods Explored" -> "Background"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Background"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "U.S. crime rates have varied, decreasing in the 1990s and rising recently in some areas. These trends differ regionally due to factors like urbanization and economic conditions. Socioeconomic issues, the opioid crisis, and drug trafficking also influence crime rates. Technological advances in law enforcement have transformed crime prevention and investigation. However, public perception of crime, often shaped by media, may not always reflect actual statistics."

# ---------------------------------------------------------------------
# 5. Slide 4: "The Methods/Tools Used" -> "Methods Explored"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Methods Explored"
Set-Paragraphs $s4.Shapes.Item(2).TextFrame.TextRange @(
    "Considered Methods: List and briefly describe the methods you considered.",
    "Method Selection: Discuss why you chose the method/tool you did."
)

# ---------------------------------------------------------------------
# 6. Slide 5: "Results" -> "The Methods/Tools Used"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "The Methods/Tools Used"
Set-Paragraphs $s5.Shapes.Item(2).TextFrame.TextRange @(
    "Method/Tool Used: Describe the chosen method/tool.",
    "Rationale: Explain the rationale behind your choice."
)

# ---------------------------------------------------------------------
# 7. Slide 6: "Lessons Learned and Future Plans" -> "Results"
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Results"
Set-Paragraphs $s6.Shapes.Item(2).TextFrame.TextRange @(
    "Preliminary Results: Present your initial findings or data.",
    "Further Results: More results or visualizations."
)
